# Scheduled-runner update: refresh Universalis price snapshots + recompute
# profit columns (H..N) for the affected Leve rows across all job sheets.
$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1661
$ws.Range("I5").Value = 1907
$ws.Range("J5").Value = 800
$ws.Range("K5").Value = 1907
$ws.Range("L5").Value = 800
$ws.Range("M5").Value = -1792
$ws.Range("N5").Value = -1030
$ws.Range("H17").Value = 3686.122
$ws.Range("J17").Value = 3686.122
$ws.Range("L17").Value = 11058.366
$ws.Range("N17").Value = -11394.366
$ws.Range("H113").Value = 2538.484
$ws.Range("I113").Value = 3042.9333
$ws.Range("J113").Value = 2065.5625
$ws.Range("K113").Value = 3042.9333
$ws.Range("L113").Value = 2065.5625
$ws.Range("M113").Value = 211.0666999999999
$ws.Range("N113").Value = -8573.5625
$ws.Range("H123").Value = 31355.334
$ws.Range("J123").Value = 31355.334
$ws.Range("L123").Value = 31355.334
$ws.Range("N123").Value = -41155.334
$ws.Range("H129").Value = 1238.4048
$ws.Range("I129").Value = 1408.2
$ws.Range("J129").Value = 1185.3438
$ws.Range("K129").Value = 4224.6
$ws.Range("L129").Value = 3556.0314
$ws.Range("M129").Value = 775.3999999999996
$ws.Range("N129").Value = -13556.0314
$ws.Range("H130").Value = 49085
$ws.Range("J130").Value = 49085
$ws.Range("L130").Value = 49085
$ws.Range("N130").Value = -59125
$ws.Range("H137").Value = 4956.657
$ws.Range("I137").Value = 1700
$ws.Range("J137").Value = 5261.9688
$ws.Range("K137").Value = 5100
$ws.Range("L137").Value = 15785.9064
$ws.Range("M137").Value = -2550
$ws.Range("N137").Value = -20885.9064
$ws.Range("H141").Value = 4999.9165
$ws.Range("I141").Value = 3555
$ws.Range("J141").Value = 9334.666999999999
$ws.Range("K141").Value = 10665
$ws.Range("L141").Value = 28004.001
$ws.Range("M141").Value = -5485
$ws.Range("N141").Value = -38364.001

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 28740
$ws.Range("J37").Value = 28740
$ws.Range("L37").Value = 28740
$ws.Range("N37").Value = -29286
$ws.Range("H61").Value = 3330.1516
$ws.Range("I61").Value = 2277.077
$ws.Range("J61").Value = 4014.65
$ws.Range("K61").Value = 2277.077
$ws.Range("L61").Value = 4014.65
$ws.Range("M61").Value = -2065.077
$ws.Range("N61").Value = -4438.65
$ws.Range("H74").Value = 1268.7858
$ws.Range("I74").Value = 646.7273
$ws.Range("K74").Value = 646.7273
$ws.Range("M74").Value = 227.2727
$ws.Range("H77").Value = 1268.7858
$ws.Range("I77").Value = 646.7273
$ws.Range("K77").Value = 3233.6365
$ws.Range("M77").Value = 1134.3635
$ws.Range("H80").Value = 53750
$ws.Range("J80").Value = 53750
$ws.Range("L80").Value = 53750
$ws.Range("N80").Value = -55746
$ws.Range("H83").Value = 53750
$ws.Range("J83").Value = 53750
$ws.Range("L83").Value = 161250
$ws.Range("N83").Value = -171234
$ws.Range("H101").Value = 46746
$ws.Range("J101").Value = 46746
$ws.Range("L101").Value = 46746
$ws.Range("N101").Value = -53236
$ws.Range("H122").Value = 1794.2174
$ws.Range("I122").Value = 1792.4
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 5377.200000000001
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -2927.200000000001
$ws.Range("N122").Value = -10300
$ws.Range("H130").Value = 43664
$ws.Range("J130").Value = 43664
$ws.Range("L130").Value = 43664
$ws.Range("N130").Value = -53704
$ws.Range("H132").Value = 13515216
$ws.Range("I132").Value = 21740340
$ws.Range("K132").Value = 65221020
$ws.Range("M132").Value = -65218490
$ws.Range("H134").Value = 51158
$ws.Range("J134").Value = 51158
$ws.Range("L134").Value = 51158
$ws.Range("N134").Value = -61298
$ws.Range("H136").Value = 3330.1516
$ws.Range("I136").Value = 2277.077
$ws.Range("J136").Value = 4014.65
$ws.Range("K136").Value = 6831.231000000001
$ws.Range("L136").Value = 12043.95
$ws.Range("M136").Value = -4281.231000000001
$ws.Range("N136").Value = -17143.95
$ws.Range("H137").Value = 52126.668
$ws.Range("J137").Value = 52126.668
$ws.Range("L137").Value = 52126.668
$ws.Range("N137").Value = -62326.668

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 33857.5
$ws.Range("J35").Value = 33857.5
$ws.Range("L35").Value = 33857.5
$ws.Range("N35").Value = -34477.5
$ws.Range("H130").Value = 40554
$ws.Range("J130").Value = 40554
$ws.Range("L130").Value = 40554
$ws.Range("N130").Value = -50594

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1831.5172
$ws.Range("I58").Value = 1520.12
$ws.Range("J58").Value = 3777.75
$ws.Range("K58").Value = 1520.12
$ws.Range("L58").Value = 3777.75
$ws.Range("M58").Value = -1317.12
$ws.Range("N58").Value = -4183.75
$ws.Range("H133").Value = 13194
$ws.Range("J133").Value = 13194
$ws.Range("L133").Value = 13194
$ws.Range("N133").Value = -18254
$ws.Range("H136").Value = 1831.5172
$ws.Range("I136").Value = 1520.12
$ws.Range("J136").Value = 3777.75
$ws.Range("K136").Value = 4560.36
$ws.Range("L136").Value = 11333.25
$ws.Range("M136").Value = -2010.36
$ws.Range("N136").Value = -16433.25
$ws.Range("H138").Value = 45115.8
$ws.Range("J138").Value = 45115.8
$ws.Range("L138").Value = 45115.8
$ws.Range("N138").Value = -55395.8
$ws.Range("H139").Value = 41563.168
$ws.Range("J139").Value = 41475.8
$ws.Range("L139").Value = 41475.8
$ws.Range("N139").Value = -51755.8

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 333333660
$ws.Range("I80").Value = 500
$ws.Range("J80").Value = 500000260
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 1500000780
$ws.Range("M80").Value = -564
$ws.Range("N80").Value = -1500002652
$ws.Range("H83").Value = 333333660
$ws.Range("I83").Value = 500
$ws.Range("J83").Value = 500000260
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 4500002340
$ws.Range("M83").Value = 180
$ws.Range("N83").Value = -4500011700
$ws.Range("H139").Value = 109927.96
$ws.Range("I139").Value = 233305.39
$ws.Range("J139").Value = 3000.8667
$ws.Range("K139").Value = 699916.17
$ws.Range("L139").Value = 9002.6001
$ws.Range("M139").Value = -694776.17
$ws.Range("N139").Value = -19282.6001

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 34628.5
$ws.Range("J110").Value = 34628.5
$ws.Range("L110").Value = 34628.5
$ws.Range("N110").Value = -42808.5
$ws.Range("H122").Value = 958.625
$ws.Range("I122").Value = 864.46155
$ws.Range("J122").Value = 1366.6666
$ws.Range("K122").Value = 2593.38465
$ws.Range("L122").Value = 4099.9998
$ws.Range("M122").Value = -143.38465
$ws.Range("N122").Value = -8999.9998
$ws.Range("H123").Value = 14998
$ws.Range("J123").Value = 14998
$ws.Range("L123").Value = 14998
$ws.Range("N123").Value = -19898
$ws.Range("H130").Value = 46566.555
$ws.Range("J130").Value = 46566.555
$ws.Range("L130").Value = 46566.555
$ws.Range("N130").Value = -56606.555
$ws.Range("H132").Value = 4525.6924
$ws.Range("I132").Value = 1652
$ws.Range("J132").Value = 7399.385
$ws.Range("K132").Value = 4956
$ws.Range("L132").Value = 22198.155
$ws.Range("M132").Value = -2426
$ws.Range("N132").Value = -27258.155

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 21259
$ws.Range("J121").Value = 21259
$ws.Range("L121").Value = 21259
$ws.Range("N121").Value = -24753
$ws.Range("H127").Value = 42224
$ws.Range("J127").Value = 42224
$ws.Range("L127").Value = 42224
$ws.Range("N127").Value = -52144
$ws.Range("H136").Value = 2500.7666
$ws.Range("I136").Value = 1920.4286
$ws.Range("J136").Value = 3854.889
$ws.Range("K136").Value = 5761.2858
$ws.Range("L136").Value = 11564.667
$ws.Range("M136").Value = -3211.2858
$ws.Range("N136").Value = -16664.667
$ws.Range("H139").Value = 83599.5
$ws.Range("J139").Value = 58132.668
$ws.Range("L139").Value = 58132.668
$ws.Range("N139").Value = -68412.66800000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1600
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1600
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1600
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -1880
$ws.Range("H86").Value = 29333
$ws.Range("J86").Value = 29333
$ws.Range("L86").Value = 29333
$ws.Range("N86").Value = -31579
$ws.Range("H89").Value = 29333
$ws.Range("J89").Value = 29333
$ws.Range("L89").Value = 146665
$ws.Range("N89").Value = -157897
$ws.Range("H94").Value = 37326
$ws.Range("J94").Value = 37326
$ws.Range("L94").Value = 37326
$ws.Range("N94").Value = -39128
$ws.Range("H128").Value = 46045.668
$ws.Range("J128").Value = 46045.668
$ws.Range("L128").Value = 46045.668
$ws.Range("N128").Value = -56005.668
$ws.Range("H132").Value = 1627.159
$ws.Range("I132").Value = 1141.4
$ws.Range("J132").Value = 2668.0715
$ws.Range("K132").Value = 3424.2
$ws.Range("L132").Value = 8004.2145
$ws.Range("M132").Value = -894.2000000000003
$ws.Range("N132").Value = -13064.2145
$ws.Range("H136").Value = 18597.69
$ws.Range("I136").Value = 33039.066
$ws.Range("J136").Value = 2016.8518
$ws.Range("K136").Value = 99117.198
$ws.Range("L136").Value = 6050.555399999999
$ws.Range("M136").Value = -96567.198
$ws.Range("N136").Value = -11150.5554
